# "Result and evaluation part of thesis"
# Applies the updates described by the diff to the first worksheet (Tabelle1):
#  - adds a new shared string "up to 30 Mb/s" and retargets H26 to it
#  - adds a new "0.0" number format + centred cell style, and applies it
#    (plus refreshed values) to the E29:H31 data block
#  - hides gridlines, and moves the sheet selection to A19
#  - widens columns E:H slightly
#  - sets the page setup (paper size / orientation) for printing

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- H26: "30 Mb/s" -> "up to 30 Mb/s" ------------------------------------
$ws.Range("H26").Value = "up to 30 Mb/s"

# --- E29:H31: new number format (0.0) + centred alignment, updated values -
$dataRng = $ws.Range("E29:H31")
$dataRng.HorizontalAlignment = -4108   # xlCenter
$dataRng.NumberFormat = "0.0"

$ws.Range("E29").Value = 1.5
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 3.5
$ws.Range("H29").Value = 3.8

$ws.Range("E30").Value = 3.9
$ws.Range("F30").Value = 4.3
$ws.Range("G30").Value = 4.7
$ws.Range("H30").Value = 4.8

$ws.Range("E31").Value = 3.7
$ws.Range("F31").Value = 4.4
$ws.Range("G31").Value = 4.6
$ws.Range("H31").Value = 4.6

# --- sheet view: hide gridlines, move selection ---------------------------
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A19").Select()

# --- widen columns E:H slightly -------------------------------------------
$ws.Range("E1:H1").ColumnWidth = 13.5703125

# --- page setup for printing ----------------------------------------------
$ws.PageSetup.PaperSize = 9       # xlPaperA4
$ws.PageSetup.Orientation = 1     # xlPortrait
